$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update row 2 values: Current_Ct_Day, Current_Pct_Ct, Current_Ct_Tokens, Current_Pct_Tokens
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 0.0002777777777777778
$ws.Range("K2").Value = 1891
$ws.Range("L2").Value = 0.003782
